$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D3").Value = 44160
$ws.Range("K3").Value = "Dina"
$ws.Range("L3").Value = "Primera"
$ws.Range("N3").Value = 20000
$ws.Range("O3").Value = 20000
$ws.Range("P3").Value = 20000
$ws.Range("Q3").Value = "$/caja 15 kilos"
$ws.Range("R3").Value = "Provincia de Limarí"
$ws.Range("S3").Value = 1333
$ws.Range("T3").Value = 15
$ws.Range("D4").Value = 44179
$ws.Range("K4").Value = "Dina"
$ws.Range("M4").Value = 150
$ws.Range("N4").Value = 18000
$ws.Range("O4").Value = 18000
$ws.Range("P4").Value = 18000
$ws.Range("Q4").Value = "$/caja 18 kilos"
$ws.Range("R4").Value = "Provincia de San Felipe de Aconcagua"
$ws.Range("S4").Value = 1000
$ws.Range("T4").Value = 18
$ws.Range("D5").Value = 44186
$ws.Range("K5").Value = "Dina"
$ws.Range("M5").Value = 150
$ws.Range("N5").Value = 15000
$ws.Range("O5").Value = 15000
$ws.Range("P5").Value = 15000
$ws.Range("Q5").Value = "$/caja 18 kilos"
$ws.Range("R5").Value = "Región Metropolitana"
$ws.Range("S5").Value = 833
$ws.Range("D6").Value = 44168
$ws.Range("M6").Value = 250
$ws.Range("N6").Value = 10000
$ws.Range("O6").Value = 10000
$ws.Range("P6").Value = 10000
$ws.Range("Q6").Value = "$/caja 10 kilos"
$ws.Range("R6").Value = "Provincia de San Felipe de Aconcagua"
$ws.Range("S6").Value = 1000
$ws.Range("T6").Value = 10
$ws.Range("D7").Value = 44168
$ws.Range("M7").Value = 100
$ws.Range("Q7").Value = "$/caja 18 kilos"
$ws.Range("S7").Value = 944
$ws.Range("T7").Value = 18
$ws.Range("L8").Value = "Primera"
$ws.Range("M8").Value = 200
$ws.Range("N8").Value = 17000
$ws.Range("O8").Value = 17000
$ws.Range("P8").Value = 17000
$ws.Range("S8").Value = 1062
$ws.Range("D9").Value = 44162
$ws.Range("K9").Value = "Castle Brite"
$ws.Range("L9").Value = "Segunda"
$ws.Range("M9").Value = 100
$ws.Range("N9").Value = 15000
$ws.Range("O9").Value = 15000
$ws.Range("P9").Value = 15000
$ws.Range("Q9").Value = "$/caja 16 kilos granel"
$ws.Range("S9").Value = 938
$ws.Range("T9").Value = 16
$ws.Range("D10").Value = 44167
$ws.Range("K10").Value = "Castle Brite"
$ws.Range("M10").Value = 300
$ws.Range("Q10").Value = "$/caja 16 kilos granel"
$ws.Range("R10").Value = "Provincia de Limarí"
$ws.Range("S10").Value = 938
$ws.Range("T10").Value = 16
$ws.Range("D11").Value = 44174
$ws.Range("M11").Value = 200
$ws.Range("D12").Value = 44189
$ws.Range("K12").Value = "Dina"
$ws.Range("M12").Value = 50
$ws.Range("N12").Value = 15000
$ws.Range("O12").Value = 15000
$ws.Range("P12").Value = 15000
$ws.Range("Q12").Value = "$/caja 18 kilos"
$ws.Range("S12").Value = 833
$ws.Range("T12").Value = 18
$ws.Range("D13").Value = 44172
$ws.Range("K13").Value = "Castle Brite"
$ws.Range("L13").Value = "Especial"
$ws.Range("M13").Value = 120
$ws.Range("Q13").Value = "$/caja 10 kilos"
$ws.Range("R13").Value = "Provincia de San Felipe de Aconcagua"
$ws.Range("S13").Value = 1500
$ws.Range("T13").Value = 10
$ws.Range("D14").Value = 44172
$ws.Range("N14").Value = 11000
$ws.Range("O14").Value = 11000
$ws.Range("P14").Value = 11000
$ws.Range("S14").Value = 1100
$ws.Range("D15").Value = 44176
$ws.Range("Q15").Value = "$/caja 18 kilos granel"
$ws.Range("D16").Value = 44161
$ws.Range("K16").Value = "Castle Brite"
$ws.Range("M16").Value = 150
$ws.Range("N16").Value = 20000
$ws.Range("O16").Value = 20000
$ws.Range("P16").Value = 20000
$ws.Range("Q16").Value = "$/caja 18 kilos granel"
$ws.Range("R16").Value = "Provincia de Limarí"
$ws.Range("S16").Value = 1111
$ws.Range("D17").Value = 44181
$ws.Range("M17").Value = 220
$ws.Range("N17").Value = 17000
$ws.Range("O17").Value = 17000
$ws.Range("P17").Value = 17000
$ws.Range("S17").Value = 944
